$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting the existing item rows (and the
# totals / footer rows below them) down by one.
$ws.Rows("8:8").Insert()

# Copy the formatting (styles) of the row below onto the freshly inserted,
# still-blank row 8 so it matches the other item rows exactly.
$ws.Range("A9:Q9").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Recreate the merged-cell layout used by every item row for the new row 8.
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# Populate row 8 with the new item: IVY ZAD SYRUP 120 ML.
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "IVY ZAD SYRUP 120 ML"
$ws.Range("H8").Value = "1:0"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "65.00"
$ws.Range("P8").Value = "65.0000"
$ws.Range("Q8").Value = "1:0"

# Renumber the "م" (item index) column for the two rows pushed down.
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4

# Restore the row heights to match the report's layout exactly.
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 25.5
$ws.Rows("12:12").RowHeight = 16.5

# Update the running total (now on row 11) to include the new item's price.
$ws.Range("P11").Value = 99.67

# Update the "generated on" timestamp in the footer (now on row 12).
$ws.Range("A12").Value = "Tuesday, 7 October, 2025 9:22 AM"
